$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the feature-engineering image input value for the existing
# experiment row (Experiment Count = 0): "Features" is now "None".
$ws.Range("C3").Value = "None"

# Add the new experiment row (Experiment Count = 1) capturing the
# re-run with the corrected image input.
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "Color (RGB) + GLCM"
$ws.Range("D4").Value = "1000/ class"
$ws.Range("E4").Value = "Yes"
$ws.Range("F4").Value = "None"
$ws.Range("G4").Value = 0.74274274274274199
$ws.Range("H4").Value = "SVC(C=1000, gamma=0.001)"
$ws.Range("I4").Value = "48m 34.6s"

# Grow Table1 so the new row (plus the following blank row) is included.
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("B2:I5"))

# Match the author's final selection.
$ws.Range("I5").Select()
